$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'فیبرونکتین در جنین نشان دهنده چیست؟'
$ws.Range('A3').Value = 'یک بیمار 62 ساله با هماچوری مراجعه می‌کند. در بررسی‌ها، سنگ کلیه، کلسیفیکاسیون در دیواره مثانه و مثانه کوچک منقبض شده مشاهده می‌شود. محتمل‌ترین علت چیست؟'
$ws.Range('A4').Value = 'داروی مورد استفاده برای تسکین درد در جراحی سرپایی (opd) کدام است؟'
$ws.Range('A5').Value = 'کدام یک از شل‌کننده‌های عضلات صاف زیر با تأثیر بر آزادسازی کلسیم عمل نمی‌کنند؟'
$ws.Range('A6').Value = 'دوكاسانول به عنوان چه موردي استفاده مي‌شود؟'
$ws.Range('A7').Value = 'داروی انتخابی در بیماری ریوی بینابینی کدام است؟'
$ws.Range('A8').Value = 'اثر آبشار در کدام ناحیه دیده می‌شود؟'
$ws.Range('A9').Value = 'افت فشار خون در آسیب حاد نخاعی به چه دلیل رخ می‌دهد؟'
$ws.Range('A10').Value = 'با وجود نیمه عمر کوتاه (2 ساعت)، مهارکننده‌های پمپ پروتون (ppis) باعث سرکوب طولانی مدت ترشح اسید (تا 48 ساعت) می‌شوند زیرا:'
$ws.Range('A11').Value = 'بیشترین جریان خون به ازای هر 100 گرم از اندام در کدام مشاهده می‌شود؟'
$ws.Range('A12').Value = 'استئوبلاست‌هایی که سطح پریودنتال استخوان آلوئولار را می‌پوشانند، تشکیل می‌دهند:'
$ws.Range('A13').Value = 'استروئیدها در کدام مورد منع مصرف دارند؟'
$ws.Range('A14').Value = 'اثر تاموکسیفن:'
$ws.Range('A15').Value = 'اولین شواهد کلسیفیکاسیون دندان پیشین جانبی فک پایین در چه زمانی رخ می‌دهد؟'
$ws.Range('A16').Value = 'dna نشاندار رادیواکتیو در یک محیط غیر رادیواکتیو دو بار تکثیر شد. کدام یک از گزینه های زیر صحیح است؟'
$ws.Range('A17').Value = 'ظرفیت کل ریه به چه چیزی بستگی دارد؟'
$ws.Range('A18').Value = 'تمام موارد زیر ممکن است در یک بارداری طبیعی مشاهده شوند به جز:'
$ws.Range('A19').Value = 'قانون ایمبه-فیک مرتبط با کدام مورد است؟'
$ws.Range('A20').Value = 'کتوزیس ناشی از چیست؟'
$ws.Range('A21').Value = 'شباهت بین fnh و آدنوم کبدی شامل همه موارد زیر است به جز'
$ws.Range('A22').Value = 'آپاندیس بیضه باقیمانده کدام است؟'
$ws.Range('A23').Value = 'فنولداپام چیست؟'
$ws.Range('A24').Value = 'سینوس مزمن ترشح کننده با ذرات استخوان در کدام مورد مشاهده می‌شود؟'
$ws.Range('A25').Value = 'میلین با کدام روش رنگ آمیزی می‌شود؟'
$ws.Range('A26').Value = 'شاخه های تنه سلیاک همه موارد زیر هستند به جز -'
$ws.Range('A27').Value = 'نمونه ای از انتقال دهنده شیمیایی پپتیدی در دستگاه گوارش که به صورت پاراکرین ترشح می شود کدام است؟'
$ws.Range('A28').Value = 'دوره کمون به تعیین همه موارد زیر کمک می‌کند به جز-'
$ws.Range('A29').Value = 'استفاده از میفپریستون می تواند تا چند هفته از بارداری برای القای سقط درمانی مؤثر باشد؟'
$ws.Range('A30').Value = 'کدام یک از موارد زیر که برای جذب آهن ضروری است، در کبد سنتز شده و از طریق ترشح صفرا به روده منتقل می‌شود؟'
$ws.Range('A31').Value = 'پوشش اپیتلیال بطن‌های مغز چیست؟'
$ws.Range('A32').Value = 'در کدام یک از موارد زیر عصب بین استخوانی خلفی آسیب می بیند؟'
$ws.Range('A33').Value = 'بیمار با هایپربیلی روبینمی غیرکونژوگه و سطح بالای اوروبیلینوژن در ادرار مراجعه می‌کند. محتمل‌ترین تشخیص چیست؟'
$ws.Range('A34').Value = 'آسیب به مجرای ادرار مردانه در زیر غشاء پرینه باعث تجمع ادرار در کدام قسمت پرینه می‌شود؟'
$ws.Range('A35').Value = 'کدام یک از موارد زیر به تنهایی در بیمار مبتلا به فئوکروموسیتوما تجویز نمی‌شود؟'
$ws.Range('A36').Value = 'رادیوگرافی یک کودک ۴ ساله نشان‌دهنده ناحیه کوچکی از شفافیت در ناحیه دو شاخه شدن مولر دوم شیری پایین است. در هنگام باز کردن کانال ریشه، خونریزی مشاهده می‌شود. محتمل‌ترین درمان کدام است؟'
$ws.Range('A37').Value = 'تومور گلوموس از کدام لایه جنینی منشأ می‌گیرد؟ سپتامبر 2004'
$ws.Range('A38').Value = 'درمان انتخابی برای هیپرآلدوسترونیسم ثانویه کدام است؟'
$ws.Range('A39').Value = 'علت ایجاد سنگ کلیه در رزکسیون گسترده روده کوچک چیست؟'
$ws.Range('A40').Value = 'روش تشخیصی انتخابی برای برونشکتازی چیست؟'
$ws.Range('A41').Value = 'شاخص فقر انسانی شامل همه موارد زیر به جز کدام است؟'
$ws.Range('A42').Value = 'مزایای بیماری تمدید شده برای چند روز ارائه می‌شود؟'
$ws.Range('A43').Value = 'انقباض عروق ریوی ناشی از هیپوکسی:'
$ws.Range('A44').Value = 'چرخش مهره‌ها در اسکولیوز در کدام حالت بررسی می‌شود؟'
$ws.Range('A45').Value = 'کرم بالغ اکینوکوکوس گرانولوزوس معمولاً در کدام یک از حیوانات زیر دیده می‌شود؟'
$ws.Range('A46').Value = 'ارنست کلین عادات را به چه دسته‌هایی تقسیم کرده است؟'
$ws.Range('A47').Value = 'در انسداد روده، بررسی‌های مورد نیاز کدامند؟'
$ws.Range('A48').Value = 'آنزیم کاتالاز در کدام قسمت یافت می‌شود؟'
$ws.Range('A49').Value = 'تغییر رنگ دندان در درمان با تتراسایکلین به دلیل تشکیل کدام ماده است؟'
$ws.Range('A50').Value = 'شایع‌ترین محل برای انجام استریلیزاسیون لوله‌ای کدام است؟'
$ws.Range('A51').Value = 'در mvp، انتظار دارید که صدای اِجکشن کلیک بیشتر برجسته شود و سوفل به صدای اول قلب نزدیک‌تر شود در بیمار'
$ws.Range('A52').Value = 'گلومرولوپاتی تغییرات مینیمال ممکن است با همه موارد زیر همراه باشد به جز -'
$ws.Range('A53').Value = 'ظاهر سنگ مجرای صفراوی مشترک (cbd) در کولانژیوگرافی با اشعه ایکس چگونه است؟'
$ws.Range('A54').Value = 'پتاسیم (k+) فراوان‌ترین کاتیون درون‌سلولی است. بالاترین غلظت پتاسیم در کدام‌یک از موارد زیر دیده می‌شود؟'
$ws.Range('A55').Value = 'خستگی شنوایی در چه فرکانسی رخ می‌دهد؟'
$ws.Range('A56').Value = 'مزایای هیسترکتومی در مول هیداتی فرم عبارتند از:'
$ws.Range('A57').Value = 'کیست واقع در محل اتصال زائده بینی میانی، زائده بینی جانبی و زائده فک بالایی کدام است؟'
$ws.Range('A58').Value = 'کدام یک از موارد زیر با التهاب ملتحمه خونریزی دهنده حاد مرتبط است؟'
$ws.Range('A59').Value = 'شما در حال انجام اولین زایمان یک نوزاد هستید. کدام یک از تغییرات گردش خون زیر معمولاً در پنج دقیقه اول پس از تولد در نوزاد رخ می‌دهد؟'
$ws.Range('A60').Value = 'یک زن ۳۰ ساله نخست‌زا با فشار خون ۱۶۰/۹۶ میلی‌متر جیوه، bun 14 mg/dl و کراتینین 1.0 mg/dl. آزمایش ادرار پروتئینوری نشان نمی‌دهد. آزمایش کامل خون و تست‌های عملکرد کبدی هیچ ناهنجاری نشان نمی‌دهند. بهترین داروی اولیه برای این بیمار کدام است؟'
$ws.Range('A61').Value = 'شایع‌ترین علت ناشنوایی در بیمار مبتلا به کارسینوم نازوفارنکس چیست؟'
$ws.Range('A62').Value = 'کلیه حدود …………… درصد از برون ده قلبی را دریافت می کند -'
$ws.Range('A63').Value = 'عوارض جانبی شبیه به سندرم اکستراپیرامیدال در کدام مورد مشاهده می‌شود؟'
$ws.Range('A64').Value = 'در مورد حلقه گوانین، همه موارد زیر صحیح هستند به جز کدام؟'
$ws.Range('A65').Value = 'در مورد لنفوم مرتبط با مخاط کدام گزینه صحیح است؟  
الف) اچ. پیلوری عامل مستعد کننده است  
ب) به شیمی درمانی حساس است  
ج) چندین لنفوم  
د) پولیپ استرومال'
$ws.Range('A66').Value = 'کدام نوع از کارسینومای پستان دوطرفه و چندکانونی است؟'
$ws.Range('A67').Value = 'در مورد سندرم وبر همه موارد زیر صحیح هستند به جز'
$ws.Range('A68').Value = 'یک مرکز فرعی به نیازهای چه جمعیتی پاسخ می‌دهد؟'
$ws.Range('A69').Value = 'در بیشتر واکنش‌ها و تعاملات سلولی، آب به دلیل ویژگی‌های زیر به عنوان یک محیط عالی شناخته می‌شود به جز:'
$ws.Range('A70').Value = 'جفت خارج کوریونی خطر کدام یک را افزایش می دهد؟'
$ws.Range('A71').Value = 'یک دختر پنج ساله با تب و التهاب ملتحمه مراجعه می‌کند. معاینه فیزیکی نشان‌دهنده اریتم دهانی و ترک‌خوردگی همراه با بثورات ماکولوپاپولر منتشر و لنفادنوپاتی گردنی است. محتمل‌ترین تشخیص چیست؟'
$ws.Range('A72').Value = 'درونی‌سازی مواد و هضم لیزوزومی پس از آن:'
$ws.Range('A73').Value = 'سوراخ ماجندی، مایع مغزی-نخاعی را از کدام ساختار میانی تخلیه میکند و از آن خارج میشود؟'
$ws.Range('A74').Value = 'شایع‌ترین عارضه تراکئوستومی کدام است؟'
$ws.Range('A75').Value = 'بیمار با اسهال، بی‌اشتهایی و سوءجذب مراجعه کرده است. نمونه‌برداری از دوازدهه او انجام شد که هایپرپلازی کریپت، آتروفی ویلی و نفوذ سلول‌های t cd8+ در اپیتلیوم را نشان داد. تشخیص احتمالی بیمار چیست؟'
$ws.Range('A76').Value = 'واکسن‌های خشک‌شده منجمد همه موارد زیر هستند به جز -'
$ws.Range('A77').Value = 'خروجی مخچه تنها از کدام سلول ها صورت می گیرد؟'
$ws.Range('A78').Value = 'خطای نوع دوم، پذیرش فرض صفر به عنوان درست است زمانی که:'
$ws.Range('A79').Value = 'یک بیمار مبتلا به فلج بل در روز سوم مراجعه می‌کند. درمانی که به او داده می‌شود عبارت است از:'
$ws.Range('A80').Value = 'کدام یک از موارد زیر تیروئیدیت خاموش است؟'
$ws.Range('A81').Value = 'پرستار چگونه می‌تواند نرمال سالین را از 10% دکستروز (در هر 100 میلی‌لیتر) آماده کند؟'
$ws.Range('A82').Value = 'اندازه تخمک چقدر است؟'
$ws.Range('A83').Value = 'کدام یک کوتاه ترین نیمه عمر را دارد؟'
$ws.Range('A84').Value = 'یک بیمار مبتلا به کلانژیت اسکلروزان اولیه دچار کلانژیوکارسینوما شده است. شایع‌ترین محل بروز کلانژیوکارسینوما کدام است؟'
$ws.Range('A85').Value = 'آزمایش رسوب ژل الک برای تشخیص سم کدام یک از ارگانیسم های زیر استفاده می شود؟'
$ws.Range('A86').Value = 'صدای بلند ناگهانی بیشتر از صداهای طولانی مدت احتمال آسیب به حلزون گوش را دارد زیرا'
$ws.Range('A87').Value = 'تمام موارد زیر در مورد نوع 1 hla صحیح است به جز؟'
$ws.Range('A88').Value = 'مردی با علائم لرز، تب و سردرد مشکوک به پنومونی "آتیپیک" است. تاریخچه نشان می‌دهد که او مرغ‌داری دارد و حدود ۲ هفته پیش تعداد زیادی از مرغ‌هایش را به دلیل یک بیماری تشخیص‌داده‌نشده از دست داده است. کدام یک از گزینه‌های زیر محتمل‌ترین تشخیص برای وضعیت این مرد است؟'
$ws.Range('A89').Value = 'متابولیسم بیحس کننده های آمیدی در کدام یک از شرایط زیر کند می شود؟'
$ws.Range('A90').Value = 'پاسخ به آهن در کم خونی فقر آهن با کدام مورد مشخص می‌شود؟'
$ws.Range('A91').Value = 'تنگی مجرای کدام یک از عروق زیر جریان خون از طریق شریان های کلیوی را به خطر می اندازد؟'
$ws.Range('A92').Value = 'تشخیص این بیمار که پس از ورزش شدید با اختلال بینایی دردناک ناگهانی مراجعه کرده است، چیست؟'
$ws.Range('A93').Value = 'po2 بیش از 8 کیلوپاسکال (اشباع هموگلوبین) چند درصد است؟'
$ws.Range('A94').Value = 'سیتوژنتیک کارسینومای سلول کلیوی کروموفیل با کدام مورد مشخص می‌شود؟'
$ws.Range('A95').Value = 'شایع ترین نئوپلاسم بدخیم در دوران نوزادی کدام است؟'
$ws.Range('A96').Value = 'حساس ترین شاخص برای عملکرد توبول های کلیوی کدام است؟'
$ws.Range('A97').Value = 'یک پسر 6 ساله با توده قابل لمس شکمی در ناحیه اپیگاستر مراجعه می‌کند. تشخیص بالینی چیست؟ (استفراغ بدون صفرا وجود دارد)-'
$ws.Range('A98').Value = 'یک کودک ۸ ساله سابقه آدنوئیدکتومی را دارد که یک سال پیش انجام شده است. او اکنون از تجمع مایع در گوش میانی شکایت دارد. کدام یک از موارد زیر محتمل‌ترین تشخیص است؟'
$ws.Range('A99').Value = 'مک کئون در قرن نوزدهم کاهش بروز بیماری‌های عفونی مانند سل را مطالعه کرد و رابطه بین کاهش بروز بیماری‌های عفونی را توضیح داد و بیان کرد که این رابطه بهتر است از نظر چه عاملی درک شود؟'
$ws.Range('A100').Value = 'اجسام آستروئیدی ممکن است توسط کدام یک تولید شوند؟'
$ws.Range('A101').Value = 'شایع‌ترین عفونت کرمی در ایدز کدام است؟'
$ws.Range('A102').Value = 'ارائه خدمات بهداشتی اولیه (phc) توسط کدام کمیته انجام شد؟'
$ws.Range('A103').Value = 'یک زن با سردرد شدید ناگهانی مراجعه می‌کند. در سی‌تی‌اسکن تشخیص خونریزی زیر عنکبوتیه داده می‌شود. شایع‌ترین محل خونریزی زیر عنکبوتیه کدام است؟'
$ws.Range('A104').Value = 'راه رفتن دست-زانو در کدام یک از بیماران دیده می‌شود؟ مارس 2013 (b, e)'
$ws.Range('A105').Value = 'واکنش به محرکی که معمولاً دردناک نیست:'
$ws.Range('A106').Value = 'فیمبرکتومی در کدام یک از روش‌های زیر انجام می‌شود؟'
$ws.Range('A107').Value = 'در مورد ماکولا لوتئا کدام گزینه صحیح است؟'
$ws.Range('A108').Value = 'کدام یک از موارد زیر می‌تواند باعث آلوپسی هم اسکاردار و هم غیراسکاردار شود؟'
$ws.Range('A109').Value = 'کدام یک از ویژگی‌های زیر تفاوت بین دلیریوم و دمانس را نشان می‌دهد؟'
$ws.Range('A110').Value = 'فونتانل قدامی در چه زمانی بسته می‌شود؟ ap 07'
$ws.Range('A111').Value = 'در هموکروماتوز، آهن در کدام یک از موارد زیر رسوب نمی‌کند؟'
$ws.Range('A112').Value = 'پاپول های گاتورون در کدام بیماری دیده می شوند؟'
$ws.Range('A113').Value = 'آسیب به عصب ناشی از c5 در شبکه بازویی منجر به چه می‌شود؟'
$ws.Range('A114').Value = 'درماتوم انگشت وسط مربوط به کدام ریشه عصبی است؟'
$ws.Range('A115').Value = 'در درمان کارسینومای پاپیلاری تیروئید، رادیو ید عمدتاً سلول‌های نئوپلاستیک را از طریق چه چیزی از بین می‌برد؟'
$ws.Range('A116').Value = 'مردمک گشاد شده در همه موارد زیر دیده می‌شود به جز'
$ws.Range('A117').Value = 'کدام گزینه در مورد هپاتوبلاستوما صحیح نیست؟'
$ws.Range('A118').Value = 'مهم‌ترین دارو برای باکتری‌های مقاوم نهفته سل'
$ws.Range('A119').Value = 'بسته شدن زودرس کدام یک از درزهای زیر می‌تواند باعث دولیکوسفالی شود؟'
$ws.Range('A120').Value = 'کدام یک از داروهای ضد آریتمی زیر در بیماران مبتلا به بیماری بینابینی ریه منع مصرف دارد؟'
$ws.Range('A121').Value = 'باز شدن سینوس اتموئید خلفی در -'
$ws.Range('A122').Value = 'adhd در دوران کودکی می‌تواند در آینده منجر به چه چیزی شود؟'
$ws.Range('A123').Value = 'عصب درگیر در هایپراکوزیس'
$ws.Range('A124').Value = 'کدام یک از اصطلاحات زیر می‌تواند برای آسیب مغزی استفاده شود؟'
$ws.Range('A125').Value = 'کدام یک از موارد زیر از ویژگی‌های پورپورای هنوخ-شونلاین (hsp) نیستند؟'
$ws.Range('A126').Value = 'هموروئید در کدام قسمت رخ می‌دهد؟'
$ws.Range('A127').Value = 'فک فوسی ناشی از چیست؟'
$ws.Range('A128').Value = 'عبارت صحیح درباره کلسیم:'
$ws.Range('A129').Value = 'بهترین شاخص برای ارزیابی سوءتغذیه مزمن کدام است؟'
$ws.Range('A130').Value = 'کدام یک از مواد زیر در خون توسط آزمایش کونکل نشان داده می‌شود؟'
$ws.Range('A131').Value = 'یک بیمار با سنکوپ جدید، فشار خون 110/95 میلی‌متر جیوه و سوفل سیستولیک خشن در پایه قلب دارد که به هر دو شریان کاروتید انتشار می‌یابد. در سمع صدای دوم قلب در پایه، کدام یک از یافته‌های زیر ممکن است مشاهده شود؟'
$ws.Range('A132').Value = 'همه موارد زیر سرکوب کننده اشتها هستند به جز:'
$ws.Range('A133').Value = 'غدد لنفاوی که از پستان تخلیه می‌شوند و در پشت عضله پکتورالیس ماینور قرار دارند، به چه عنوانی شناخته می‌شوند؟'
$ws.Range('A134').Value = 'همه موارد زیر در مورد سندرم گرادنیگو صحیح است به جز؟'
$ws.Range('A135').Value = 'فعالیت ویژه یک آنزیم در کدام یک از واحدهای اندازه‌گیری زیر گزارش می‌شود؟'
$ws.Range('A136').Value = 'تست فرستنبرگ در کدام مورد مثبت دیده می‌شود؟'
$ws.Range('A137').Value = 'ضایعه کاپوسی واریسیلیفورم در کدام یک از موارد زیر دیده می‌شود؟ الف) درماتیت آتوپیک ب) بیماری داریه ج) لیکن پلان د) واریسلا زوستر'
$ws.Range('A138').Value = 'کدام عضله چرخاننده داخلی چشم است؟'
$ws.Range('A139').Value = 'یک پسر 5 ساله با سابقه تشنج های کانونی که با دوره های اختلال هوشیاری همراه است، مراجعه کرده است. بیمار همچنین سابقه تغییرات رفتاری دوره ای شامل تحریک پذیری و حالت رویا مانند همراه با حرکات شدید اندام ها را دارد. کودک این دوره ها را به خاطر نمی آورد و عامل محرک واضحی وجود ندارد. همچنین سابقه تشنج در خانواده وجود دارد. تشنج ها اغلب مقاوم به درمان هستند. بررسی های آزمایشگاهی نشان دهنده اسپایک های یک طرفه یا دو طرفه تمپورال قدامی در eeg، هیپومتابولیسم در pet بین حمله ای، هیپوپرفیوژن در spect بین حمله ای و نقایص حافظه وابسته به ماده در تست آموباربیتال داخل جمجمه ای (تست وادا) بود. mri مغز نیز انجام شد. کدام یک از ساختارهای زیر معمولاً در این شرایط درگیر می شود؟'
$ws.Range('A140').Value = 'فسفاتاز قلیایی لکوسیتی در همه موارد زیر افزایش می یابد، به جز-'
$ws.Range('A141').Value = 'در مورد حلقه اینگوینال عمیق کدام گزینه صحیح است؟'
$ws.Range('A142').Value = 'عوارض جانبی ssri'
$ws.Range('A143').Value = 'همه موارد زیر از ویژگی‌های اگزتروفی مثانه هستند به جز:'
$ws.Range('A144').Value = 'دوره نهفتگی بیماری اوریون چقدر است؟'
$ws.Range('A145').Value = 'جایگزین اسید آمینه برای تریپتوفان بدون تغییر ویژگی‌های پروتئین کدام است؟'
$ws.Range('A146').Value = 'گلیوبلاستوما مولتیفرم در موارد زیر ممکن است رخ دهد به جز:'
$ws.Range('A147').Value = 'همه موارد زیر برای پیشگیری از بارداری در کم خونی داسی شکل استفاده می‌شوند به جز'
$ws.Range('A148').Value = 'زانوی کپسول داخلی شامل چه چیزی است؟'
$ws.Range('A149').Value = 'تست توبی آیر در کدام مورد مشاهده می‌شود؟'
$ws.Range('A150').Value = 'اریتم ندوزوم در همه موارد زیر دیده می‌شود به جز -'
$ws.Range('A151').Value = 'کلسی تونین داخل بینی در چه شرایطی تجویز می‌شود؟'
$ws.Range('A152').Value = 'سلول‌های کوپفر در کبد مسئول کدام عملکرد هستند؟ سپتامبر ۲۰۰۷'
$ws.Range('A153').Value = 'کدام یک از گزینه‌های زیر از محل‌های اصلی عفونت گنوکوکی حاد نیست؟'
$ws.Range('A154').Value = 'تصویر رادیوگرافی دیستال تیبیا و فیبولای سمت راست نشان داده شده است که احتمالاً نمایانگر کدام مورد است؟'
$ws.Range('A155').Value = 'جفت هورمونی ترشح می‌کند که در تحریک رشد مجاری در غده پستانی در دوران بارداری مشارکت دارد. این هورمون کدام است؟'
$ws.Range('A156').Value = 'میزبان شایع بالانتی دیوم کولی چیست؟'
$ws.Range('A157').Value = 'درمان خط اول بهینه برای کارسینوم سلول سنگفرشی کانال مقعد شامل چه مواردی است؟'
$ws.Range('A158').Value = 'معیارهای metcalfe برای چه موردی استفاده می‌شود؟'
$ws.Range('A159').Value = 'ریز مغذی مرتبط با راش و اسهال کدام است؟'
$ws.Range('A160').Value = 'شایع‌ترین محل بروز انفارکتوس میوکارد کدام است؟'
$ws.Range('A161').Value = 'کدام ویتامین برای انتقال واحد 1-کربن مورد نیاز است؟'
$ws.Range('A162').Value = 'هیدروسل مادرزادی چگونه درمان می‌شود؟'
$ws.Range('A163').Value = 'در مورد کمبود آلفا-1 آنتی تریپسین، کدام گزینه/گزینه‌ها صحیح است؟  
الف) اتوزومال غالب  
ب) آمفیزم ریوی  
ج) سلول‌های کبدی مقاوم به دیاستاز  
د) سلول‌های کبدی با رنگ‌آمیزی اورسین مثبت  
ه) مرتبط با آنوریسم توت‌مانند'
$ws.Range('A164').Value = 'شاخص سالیوان چیست؟'
$ws.Range('A165').Value = 'زخم‌های آتشفشانی در مری در کدام مورد دیده می‌شوند؟'
$ws.Range('A166').Value = 'نمک یددار برای جلوگیری از گواتر به چه کسانی داده می‌شود؟'
$ws.Range('A167').Value = 'سلول های بنیادی در کدام یک از مکان های زیر در بدن قرار دارند؟'
$ws.Range('A168').Value = 'فشار ساکشن نباید از ----- میلی‌متر جیوه تجاوز کند، هنگام ساکشن راه هوایی نوزادان در احیای نوزادی:'
$ws.Range('A169').Value = 'ایزوآنزیم غالب ldh در عضله قلب کدام است؟'
$ws.Range('A170').Value = 'نجات اندام را می‌توان در همه موارد زیر انجام داد به جز -'
$ws.Range('A171').Value = 'در مقایسه با کولیماسیون دایره‌ای، کولیماسیون مستطیلی میزان پرتوگیری را ↓ می‌کند به میزان:'
$ws.Range('A172').Value = 'یک اسیدوز متابولیک با گپ آنیونی طبیعی در بیماران مبتلا به کدام یک از موارد زیر رخ می‌دهد؟'
$ws.Range('A173').Value = 'کدام مورد در درمان گال استفاده نمی‌شود؟'
$ws.Range('A174').Value = 'همه موارد زیر سیتوکین های پیروژن هستند، به جز-'
$ws.Range('A175').Value = 'هاراکیری مرگ به روش:'
$ws.Range('A176').Value = 'یک سرباز ارتش که در یک منطقه جنگلی دورافتاده مستقر بود، تب و سردرد داشت. تب او 104 درجه فارنهایت و نبضش 70 ضربه در دقیقه بود. در زمان مراجعه به بیمارستان ارجاعی، یک ضایعه اریتماتوز حدود 1 سانتیمتری روی پا همراه با وزیکول‌های کوچک و لنفادنوپاتی عمومی داشت. نمونه خون او برای انجام سرولوژی جهت تشخیص بیماری ریکتزیایی جمع‌آوری شد. کدام یک از نتایج زیر در واکنش ویل-فلیکس در این شرایط بالینی تشخیصی خواهد بود؟'
$ws.Range('A177').Value = 'فاصله زمانی بین معرفی باکتری به محیط کشت و شروع تکثیر آن چیست؟'
$ws.Range('A178').Value = 'نسیدیوبلاستوما ناشی از هایپرپلازی کدام سلول است؟'
$ws.Range('A179').Value = 'کدام یک از روش‌های زیر بهترین روش برای اندازه‌گیری سطح تیامین در خون است؟'
$ws.Range('A180').Value = 'مخاط دهانه رحم با چه نوع اپیتلیومی پوشیده شده است؟'
$ws.Range('A181').Value = 'بهترین روش برای کشت m.leprae چیست؟'
$ws.Range('A182').Value = 'سولفونامید با کدام فرآیند کونژوگه می‌شود؟'
$ws.Range('A183').Value = 'کدام گزینه درباره گلیکولیز نادرست است؟'
$ws.Range('A184').Value = 'کدام آنزیم توسط سدیم فلوراید مهار می‌شود؟'
$ws.Range('A185').Value = 'همه موارد زیر در مدیاستینوم خلفی یافت می‌شوند به جز:'
$ws.Range('A186').Value = 'سم وبا ناشی از کدام مورد زیر است؟'
$ws.Range('A187').Value = 'داروی مورد استفاده برای مالاریای ویواکس مقاوم به کلروکین همه موارد زیر هستند به جز: سپتامبر 2009'
$ws.Range('A188').Value = 'یک مارپیچ آلفا در یک پروتئین به احتمال زیاد مختل می‌شود اگر یک جهش معنی‌دار اسید آمینه زیر را در ساختار مارپیچ آلفا وارد کند:'
$ws.Range('A189').Value = 'تمام موارد زیر از اثرات بتا بلاکرها در نارسایی قلبی (chf) هستند به جز:'
$ws.Range('A190').Value = 'خطوط زان در کدام یک از موارد زیر دیده می‌شوند؟'
$ws.Range('A191').Value = 'اوکرونسیس در کدام بیماری یافت می‌شود؟'
$ws.Range('A192').Value = 'انفارکتوس مغزی ناشی از چیست؟'
$ws.Range('A193').Value = 'نادرست در مورد h. influenzae:'
$ws.Range('A194').Value = 'کدام یک از موارد زیر یک سازمان داوطلبانه سلامت نیست؟'
$ws.Range('A195').Value = 'جنسیت و سن با کدام روش نمایش داده می‌شود؟'
$ws.Range('A196').Value = 'میانگین زمان مورد نیاز برای بهبود کامل انفارکتوس میوکارد چقدر است؟'
$ws.Range('A197').Value = 'آنتی‌بادی علیه hsv شروع به ظاهر شدن می‌کند'
$ws.Range('A198').Value = 'هنگامی که یک دارو به گیرنده متصل می‌شود و باعث ایجاد اثری مخالف آگونیست می‌شود، این عمل چه نامیده می‌شود؟'
$ws.Range('A199').Value = 'مری در کدام مدیاستینوم قرار دارد؟'
$ws.Range('A200').Value = 'کدام یک از موارد زیر مناسب‌ترین موقعیت برای تجویز قرص فقط پروژسترونی است؟'
$ws.Range('A201').Value = 'یک بیمار مبتلا به csom دارای کلهستئاتوم است و با سرگیجه مراجعه می‌کند. درمان انتخابی کدام است؟'
$ws.Range('A202').Value = 'بیماری روز دوشنبه در اثر مواجهه با کدام ماده ایجاد می‌شود؟'
$ws.Range('A203').Value = 'یک کودک 3 ساله کاملاً واکسینه نشده برای اولین بار به کلینیک واکسیناسیون در مرکز بهداشت اولیه مراجعه می‌کند. او باید دریافت کند -'
$ws.Range('A204').Value = 'مکانیسم اثر فلوروکینولون ها چیست؟'
$ws.Range('A205').Value = 'تومور وارتین چیست؟'
$ws.Range('A206').Value = 'همه موارد زیر از ویژگی‌های پاتولوژیک مهم در جهش ژن atp7b هستند، به جز:'
$ws.Range('A207').Value = 'کدام یک از موارد زیر را نمی‌توان با اسپیرومتر ساده اندازه‌گیری کرد؟'
$ws.Range('A208').Value = 'کدام دارو برای جلوگیری از انتقال hiv از مادر به کودک، به ویژه در دوران شیردهی، تجویز می‌شود؟'
$ws.Range('A209').Value = 'کدام یک از موارد زیر از اهداف تریتیوریشن (آماده سازی آمالگام) نیست؟'
$ws.Range('A210').Value = 'اکسیداسیون 1 مول از اسید چرب c16، پالمیتات، چند مولکول atp تولید می‌کند؟'
$ws.Range('A211').Value = 'کدام یک از عبارات زیر در مورد اتصال شکافدار (گَپ جانکشن) صحیح است؟'
$ws.Range('A212').Value = 'استفاده طولانی مدت از کدام داروی ادرارآور می‌تواند منجر به ژنیکوماستی شود؟'
$ws.Range('A213').Value = 'هنگامی که یک فرد با دیگری تماس برقرار می‌کند، انتقال مواد از یکی به دیگری رخ می‌دهد. این پدیده چه نامیده می‌شود؟'
$ws.Range('A214').Value = 'متاستاز به غدد لنفاوی گردن از کدام سرطان ناشی می‌شود؟'
$ws.Range('A215').Value = 'بیمار از درد متناوب زانو شکایت دارد که گاهی با تورم زانو همراه است. با بررسی بیشتر، سابقه سیاه شدن ادرار در معرض هوا گزارش شده است. بررسی‌های روتین نتیجه خاصی نداشته است. آنزیم درگیر ممکن است کدام یک باشد؟'
$ws.Range('A216').Value = 'مهاجرت گلبول‌های سفید عمدتاً از طریق کدام یک از موارد زیر رخ می‌دهد؟'
$ws.Range('A217').Value = 'کدام گزینه در مورد اندوکاردیت عفونی صحیح نیست؟'
$ws.Range('A218').Value = 'فعال شدن کدام یک از گیرنده های زیر باعث انقباض عروق می شود؟'
$ws.Range('A219').Value = 'کدام یک از موارد زیر در مطالعات اپیدمیولوژیک تحت تأثیر تغییرات بین مشاهده‌گران قرار می‌گیرد؟'
$ws.Range('A220').Value = 'کدام یک از معرف‌های زیر برای تعیین اسید آمینه انتهای n یک پلی‌پپتید مفیدتر است؟'
$ws.Range('A221').Value = 'کدام یک از موارد زیر در تشخیص توده لوله‌ای در حاملگی خارج رحمی مفید نیست؟'
$ws.Range('A222').Value = 'کدام یک از موارد زیر در مورد انفارکتوس قرمز صحیح نیست؟'
$ws.Range('A223').Value = 'تخم کدام یک از انگل‌های زیر دارای رشته‌های قطبی است که از دو طرف آمبروفور منشأ می‌گیرند؟'
$ws.Range('A224').Value = 'استفاده طولانی مدت از استروئیدهای موضعی باعث کدام مورد می‌شود؟'
$ws.Range('A225').Value = 'بررسی مداوم عوامل مؤثر بر وقوع یک بیماری چیست؟'
$ws.Range('A226').Value = 'بررسی زجاجیه بهتر است با چه وسیله ای انجام شود؟'
$ws.Range('A227').Value = 'کدام یک از موارد زیر شدیدترین نوع گال است؟'
$ws.Range('A228').Value = 'نکروز آرتریولیت با نکروز فیبرینوئید ناشی از کدام مکانیسم است؟'
$ws.Range('A229').Value = 'شایع‌ترین علت مرگ در بیماری منکس چیست؟'
$ws.Range('A230').Value = 'یک مرد میانسال از درد بالای شکم پس از یک وعده غذایی سنگین شکایت دارد. در معاینه، حساسیت در ناحیه بالای شکم مشاهده می‌شود و در عکس‌برداری اشعه ایکس، پهن‌شدگی مدیاستینوم همراه با پنومومدیاستینوم دیده می‌شود. تشخیص چیست؟'
$ws.Range('A231').Value = 'ارزش بیولوژیکی بیشترین مقدار را در کدام گزینه دارد؟'
$ws.Range('A232').Value = 'در پورپورای هنوخ-شونلاین، درگیری کلیوی معمولاً مشاهده نمی‌شود اگر تا چه زمانی درگیری وجود نداشته باشد؟'
$ws.Range('A233').Value = 'کدام گزینه در مورد توزیع سن و جنسیت در اسکیزوفرنی صحیح است؟'
$ws.Range('A234').Value = 'سناریوی بالینی با کیسه کرم‌مانند -'
$ws.Range('A235').Value = 'همه به سینوس کرونری باز می‌شوند به جز -'
$ws.Range('A236').Value = 'وزن بسیار کم هنگام تولد -'
$ws.Range('A237').Value = 'علامت "o" در ادیوگرام تن خالص نشان‌دهنده چیست؟'
$ws.Range('A238').Value = 'تحسین، حالات چهره مثبت، دست دادن با بیمار در کدام دسته قرار می‌گیرد؟'
$ws.Range('A239').Value = 'بر اساس سیاست ملی جمعیت سال 2000، هدف میان‌مدت کاهش نرخ باروری کلی به سطح جایگزینی تا کدام سال است؟'
$ws.Range('A240').Value = 'کدام یک از ارگانیسم های زیر مسئول بیماری "ریه کارگر مالت" است؟'
$ws.Range('A241').Value = 'گلوتارآلدئید به فرمالدئید ترجیح داده می‌شود زیرا:'
$ws.Range('A242').Value = 'مهم‌ترین عامل پیش‌آگهی برای سرطان پستان چیست؟'
$ws.Range('A243').Value = 'پس از طحالبرداری، شایعترین عفونت کدام است؟'
$ws.Range('A244').Value = 'ساختار اصلی که حمایت کننده رحم است کدام است؟'
$ws.Range('A245').Value = 'کدام یک از موارد زیر جزو بیماری‌های وابسته به x نیست؟  
(a) سندرم بی‌حساسیتی به آندروژن  
(b) دیستروفی عضلانی دوشن  
(c) هموفیلی  
(d) بتا تالاسمی  
(e) بیماری ویلسون'
$ws.Range('A246').Value = 'شایع‌ترین محل هموژیوم بینی:'
$ws.Range('A247').Value = 'کدام یک از عبارات زیر در مورد بیماری انسدادی ریه صحیح است؟'
$ws.Range('A248').Value = 'محتوای پروتئین کمتر از 4 میلی‌گرم در میلی‌لیتر در کدام یک از موارد زیر مشاهده می‌شود؟'
$ws.Range('A249').Value = 'کدام یک از اختلالات ژنتیکی زیر در مردان شایع تر از زنان است؟'
$ws.Range('A250').Value = 'حداکثر دوز لیدوکائین همراه با اپی نفرین برای بی‌حسی نفوذی چقدر است؟'
$ws.Range('A251').Value = 'نوع انتقال درگیر در مورد طاعون در کک موش چیست؟'
$ws.Range('A252').Value = 'در یک ویزیت واحد، به یک کودک 9 ماهه واکسینه نشده می‌توان واکسن‌های زیر را تزریق کرد -'
$ws.Range('A253').Value = 'تکثیر ناپیوسته dna که در طول تکثیر رخ می‌دهد، منجر به تولید قطعات کوچک dna می‌شود که به نام زیر شناخته می‌شوند:'
$ws.Range('A254').Value = 'پروتز در سر استخوان فمور در چه موردی اعمال می‌شود؟'
$ws.Range('A255').Value = 'شایع‌ترین مفاصل درگیر در نوروپاتی شارکو در بیماران دیابتی کدامند؟'
$ws.Range('A256').Value = 'علامت نیلوفر آبی در رادیوگرافی قفسه سینه در کدام بیماری دیده می‌شود؟'
$ws.Range('A257').Value = 'اگر آسیب در ریشه c7 رخ دهد، در کدام قسمت بازو حس از بین می‌رود؟'
$ws.Range('A258').Value = 'واکسن مننژوکوکی برای همه موارد زیر موجود است، به جز -'
$ws.Range('A259').Value = 'شایع‌ترین بخش درگیر در اتواسکلروز کدام است؟'
$ws.Range('A260').Value = 'کدام یک از شرایط ارثی زیر باعث هیپربیلی روبینمی مستقیم می‌شود؟'
$ws.Range('A261').Value = 'منبع اصلی انرژی در ۱ دقیقه کدام است؟'
$ws.Range('A262').Value = 'هایپرگلیسمی توسط همه موارد زیر به جز ... ایجاد می‌شود:'
$ws.Range('A263').Value = 'کدام یک از موارد زیر مهارکننده فرروکلاتاز است؟'
$ws.Range('A264').Value = 'tgf-β در تمام فرآیندهای رگزایی نقش دارد به جز:'
$ws.Range('A265').Value = 'فردی به نام ''x'' در اثر تحریک، فرد دیگری به نام ''y'' را با چوبدستی میزند. این عمل منجر به تشکیل کبودی به اندازه 3 سانتیمتر در 3 سانتیمتر روی ساعد میشود. هیچ آسیب دیگری مشاهده نشده است. کدام یک از گزینههای زیر در مورد مجازات او صحیح است؟'
$ws.Range('A266').Value = 'کدام یک از موارد زیر جزء محتویات کانال ادکتور نیست؟'
$ws.Range('A267').Value = 'پدر فیزیولوژی کیست؟'
$ws.Range('A268').Value = 'ساختار عبوری از سوراخ مگنوم کدام است؟'
$ws.Range('A269').Value = 'تحریک تنفسی در کدام مورد مشاهده می‌شود؟'
$ws.Range('A270').Value = 'پلاکت چگونه به تثبیت لخته کمک می‌کند؟'
$ws.Range('A271').Value = 'در کدام یک از موارد زیر فرسایش مفاصل دیده نمی‌شود؟'
$ws.Range('A272').Value = 'یک بیمار جوان با تاری تدریجی دید در چشم چپ به بخش سرپایی چشم پزشکی مراجعه می‌کند. معاینه با لامپ اسلیت رسوبات ستاره‌ای ریز قرنیه‌ای، فلیر زلالیه و یک کاتاراکت زیرکپسولی خلفی پیچیده را نشان می‌دهد. هیچ سینکیای خلفی مشاهده نشده است. محتمل‌ترین تشخیص چیست؟'
$ws.Range('A273').Value = 'استافیلوکوک اورئوس مقاوم به متی سیلین چیست؟'
$ws.Range('A274').Value = 'ممنوعیت مطلق استفاده از آییودی چیست؟'
$ws.Range('A275').Value = 'اینوگرام چند ساعت پس از تولد انجام می‌شود؟'
$ws.Range('A276').Value = 'در مورد آمنوره ثانویه که پس از مصرف استروژن و پروژسترون، خونریزی قطعی رخ نمیدهد، مشکل در کدام سطح است؟'
$ws.Range('A277').Value = 'کدام یک از موارد زیر باعث کاهش ظرفیت انتشار ریوی (dl) نمی‌شود؟'
$ws.Range('A278').Value = 'هر دو نوع انفارکتوس رنگ پریده و خونریزی دهنده در کدام عضو یافت می‌شوند؟'
$ws.Range('A279').Value = 'یک کودک دختر با اندام تناسلی مبهم و هایپرپیگمانتاسیون پوست مراجعه می‌کند. او دارای فشار خون بالا و هیپوکالمی است. کمبود کدام یک از آنزیم‌های زیر در این کودک وجود دارد؟'
$ws.Range('A280').Value = 'کدام یک از موارد زیر در مورد پلاسمودیوم فالسیپاروم صحیح است؟'
$ws.Range('A281').Value = 'بزرگ‌ترین حجم ریه توسط کدام مورد تعریف می‌شود؟'
$ws.Range('A282').Value = 'یک بیمار مرد 50 ساله از کاهش توانایی باز کردن دهان شکایت دارد. بیمار سابقه تصادف رانندگی یک هفته پیش را ذکر می‌کند. یک عکس پانورامیک با استفاده از صفحه تقویت‌کننده حاوی گادولینیوم و لانتانیم گرفته شده است. در مقایسه با صفحه تنگستات کلسیم، صفحه مورد استفاده در این بیمار میزان تابش بیمار را به میزان زیر کاهش می‌دهد:'
$ws.Range('A283').Value = 'مدیولوس (آپکس) به کدام جهت اشاره دارد؟'
$ws.Range('A284').Value = 'همه موارد زیر در مورد آگاماگلوبولینمی صحیح هستند به جز -  
الف) از دست دادن مرکز جوانه‌ای در گره لنفاوی  
ب) لنفوسیت‌های قشری طبیعی  
ج) لنفوسیت‌های قشری طبیعی در پاراکورتکس و مدولا  
د) کاهش پولپ قرمز در طحال  
ه) نقص ایمنی (سلولی)'
$ws.Range('A285').Value = 'آزمایش ترجیحی برای تأیید ریشه‌کنی هلیکوباکتر پیلوری کدام است؟'
$ws.Range('A286').Value = 'تشکیل سنگ کیسه صفرا در همه موارد زیر رخ می‌دهد به جز'
$ws.Range('A287').Value = 'تومور بیضه در مرد 70 ساله به احتمال زیاد چیست؟'
$ws.Range('A288').Value = 'انسولین بر کدام آنزیم در گلیکولیز اثر می‌گذارد؟'
$ws.Range('A289').Value = 'کدام عکسبرداری با اشعه ایکس باید برای یک پسر 16 ساله انجام شود؟'
$ws.Range('A290').Value = 'داروی انتخابی در بیماران مبتلا به سندرم ولف-پارکینسون-وایت با فیبریلاسیون دهلیزی کدام است؟'
$ws.Range('A291').Value = 'محرک رفلکس کششی چیست؟'
$ws.Range('A292').Value = 'ظاهر "موی ایستاده" در کدام یک از موارد زیر دیده می‌شود؟'
$ws.Range('A293').Value = 'یک مرد 50 ساله که به عنوان آشپز در هتل کار می‌کند و چهار عضو خانواده تحت تکفل دارد، با تشخیص سرطان سلول سنگفرشی در مراحل اولیه کانال مقعد مواجه شده است. شانس بهبودی او بیش از 60٪ است. بهترین گزینه درمانی چیست؟'
$ws.Range('A294').Value = 'یک زن 26 ساله با شکایت از درد و تورم پای راست که طی 2 روز گذشته بدتر شده است، به اورژانس مراجعه می‌کند. او 2 روز پیش پس از یک سفر پیاده‌روی، 8 ساعت با ماشین رانندگی کرده و سپس متوجه درد در پا شده است. ابتدا فکر می‌کرد درد ناشی از فعالیت بدنی است، اما درد در طول روز بدتر شده است. سابقه پزشکی او شامل مشکل در باردار شدن و دو سقط خودبه‌خودی قبلی است. معاینه فیزیکی او نشان‌دهنده علائم حیاتی طبیعی، معاینه قلب و ریه طبیعی است. پای راست او از وسط ران به پایین متورم و حساس است. مطالعات داپلر یک ترومبوز وریدی عمیق بزرگ در وریدهای فمورال و ایلیاک که به لگن گسترش یافته است را نشان می‌دهد. آزمایش‌ها نشان‌دهنده الکترولیت‌های طبیعی، شمارش طبیعی گلبول‌های سفید و پلاکت‌ها، pt طبیعی و aptt 3 برابر طبیعی است. تست بارداری او منفی است. درمان با هپارین با وزن مولکولی کم در اورژانس آغاز می‌شود. درمان بعدی باید شامل کدام یک از موارد زیر باشد؟'
$ws.Range('A295').Value = 'تومور وارتین چیست؟'
$ws.Range('A296').Value = 'در نفریت بینابینی حاد، کدام پروتئین مرتبط است؟'
$ws.Range('A297').Value = 'کدام گزینه در مورد سوپرآنتی‌ژن صحیح نیست؟'
$ws.Range('A298').Value = 'در "شاخص تایلستروپ فیجرسکوف"، زمانی که "خطوط سفید کدر بیشتر مشخص شده و اغلب ادغام می‌شوند تا مناطق کوچک ابری شکل را در سرتاسر سطح تشکیل دهند"، نمره داده شده چیست؟'
$ws.Range('A299').Value = 'یک کودک 2 ساله با وزن 6.4 کیلوگرم و کمبود ویتامین a دارد. درجه سوءتغذیه در این کودک کدام است؟'
$ws.Range('A300').Value = 'پروتئین‌های فیبریلاری گلیال در کدام‌یک وجود دارند؟'
$ws.Range('A301').Value = 'در سوختگی درجه سه همه موارد زیر دیده می‌شوند به جز'
$ws.Range('A302').Value = 'همه موارد زیر از ویژگی‌های سندرم گیلن باره هستند، به جز:'
$ws.Range('A303').Value = 'عصب پتروزال عمیق از کجا تشکیل شده است؟'
$ws.Range('A304').Value = 'کدام یک از مواد ضدعفونی کننده زیر در از بین بردن اسپورها مؤثر نیست؟'
$ws.Range('A305').Value = 'غده مول اصلاح شده کدام غده است؟'
$ws.Range('A306').Value = 'کلسیفیکاسیون پاپ کورن در کدام مورد دیده می‌شود؟'
$ws.Range('A307').Value = 'غلظت مولی محلول خوراکی جبران مایعات چیست؟'
$ws.Range('A308').Value = 'کدام یک محیط غنی‌سازی نیست؟'
$ws.Range('A309').Value = 'رنگ کیت 3 برای بیماری‌های مقاربتی تحت برنامه کنترل ایدز چیست؟'
$ws.Range('A310').Value = 'بیماری که از سی‌تی اسکن زیر تشخیص داده می‌شود، پس از حمله کدام یک از موارد زیر رخ می‌دهد؟'
$ws.Range('A311').Value = 'بهترین روش تشخیصی برای بیماری بِست چیست؟'
$ws.Range('A312').Value = 'کدام یک از عبارات زیر که اندونوکلئازهای محدودکننده را توصیف می‌کنند صحیح است؟'
$ws.Range('A313').Value = 'ندول های ماهواره ای در کدام یک از موارد زیر دیده می شوند؟'
$ws.Range('A314').Value = 'همه در انتقال سریع آکسونی نقش دارند، به جز'
$ws.Range('A315').Value = 'عارضه ناشی از ضربه به ناحیه خطرناک صورت در میان گزینه‌های زیر کدام است؟'
$ws.Range('A316').Value = 'عبارت صحیح در مورد مایع خارج سلولی به مایع داخل سلولی کدام است؟'
$ws.Range('A317').Value = 'شایع‌ترین عضو درگیر در فیبروز رتروپریتوان چیست؟'
$ws.Range('A318').Value = 'سندرم واترهاوس-فریدریکسن در کدام یک دیده می‌شود؟'
$ws.Range('A319').Value = 'اجسام میکولیچ و راسل مشخصه کدام بیماری هستند؟'
$ws.Range('A320').Value = 'کدام یک از بیماری‌های مادرزادی قلب (chd) منجر به شانت چپ به راست می‌شود که معمولاً با سیانوز همراه است؟'
$ws.Range('A321').Value = 'در لاواژ آنترال، دسترسی از طریق کدام قسمت انجام می‌شود؟'
$ws.Range('A322').Value = 'همه موارد زیر در مسمومیت با دیژیتال دیده می‌شوند به جز'
$ws.Range('A323').Value = 'کریستال های اگزالات در ادرار در کدام مورد دیده می شود؟ (pgi dec 2006)'
$ws.Range('A324').Value = 'سندرم کروپ معمولاً توسط چه عامل ایجاد میشود؟  
الف) راینوویروس‌ها  
ب) ویروس کوکساکی a  
ج) ویروس کوکساکی b  
د) پاراآنفلوانزا'
$ws.Range('A325').Value = 'عامل ایجاد کننده گرانولوم ونروم کدام است؟'
$ws.Range('A326').Value = 'همه موارد زیر در مورد سفوروکسیم صحیح است به جز:'
$ws.Range('A327').Value = 'تمام موارد زیر در مورد آنوریسم آئورت صحیح است به جز:'
$ws.Range('A328').Value = 'بهترین آزمایش برای تشخیص سوءجذب کربوهیدرات‌ها در روده کوچک کدام است؟'
$ws.Range('A329').Value = 'کدام دارو عوارض هماتولوژیک ایجاد می‌کند؟'
$ws.Range('A330').Value = 'انتقال گروه آمینو از یک اسید آمینه به یک آلفا کتو اسید توسط کدام آنزیم انجام می‌شود؟'
$ws.Range('A331').Value = 'تمامی عبارات زیر در مورد ملیوئیدوز صحیح هستند، به جز:'
$ws.Range('A332').Value = 'مننژوانسفالیت ائوزینوفیلیک توسط کدام عامل ایجاد می‌شود؟'
$ws.Range('A333').Value = 'مورفولوژی all l3 یک بدخیمی است که از کدام رده سلولی نشات می گیرد؟'
$ws.Range('A334').Value = 'روشی که در آن نمونه از هر لایه از پیش تعریف شده جامعه گرفته می‌شود، چه نام دارد؟'
$ws.Range('A335').Value = 'تمام موارد زیر در مورد اسپوندیلیت آنکیلوزان صحیح است به جز:'
$ws.Range('A336').Value = 'در سنتز اسید چرب، از دست دادن co2 در کدام مرحله رخ می‌دهد؟'
$ws.Range('A337').Value = 'همه موارد زیر می‌توانند منجر به ژنیکوماستی شوند به جز؛'
$ws.Range('A338').Value = 'شایع ترین اختلال متابولیک در سیروز کبدی کدام است؟'
$ws.Range('A339').Value = 'دررفتگی شایع آرنج'
$ws.Range('A340').Value = 'همه موارد زیر از علل نارسایی تنفسی هایپرکاپنیک هستند به جز'
$ws.Range('A341').Value = 'اگر تست هم اکولت برای غربالگری سرطان کولون منفی باشد، هیچ آزمایش بیشتری انجام نمی‌شود. اگر تست هم اکولت مثبت باشد، فرد یک نمونه مدفوع دوم را با تست هم اکولت ii آزمایش می‌کند. اگر این نمونه دوم مثبت باشد، فرد برای ارزیابی گسترده‌تر ارجاع داده می‌شود. تأثیر حساسیت خالص و ویژگی خالص این روش غربالگری چیست؟'
$ws.Range('A342').Value = 'حداقل عمق برای پوشش چاه بهداشتی چقدر است؟'
$ws.Range('A343').Value = 'کدام یک از تکنیک‌های مسواک زدن زیر بیشترین احتمال را برای تمیز کردن شیار لثه دارد؟'
$ws.Range('A344').Value = 'طناب اسپرماتیک شامل همه موارد زیر به جز کدام است؟'
$ws.Range('A345').Value = 'بستر معده توسط همه موارد زیر تشکیل شده است به جز؟'
$ws.Range('A346').Value = 'کدام گزینه در مورد نورومیوپاتی بیماری بحرانی صحیح است؟ (دو بار در امتحان پرسیده شده است)'
$ws.Range('A347').Value = 'آنومالی ابشتاین ممکن است به عنوان عارضه کدام دارو در صورت استفاده در سه ماهه اول بارداری رخ دهد:'
$ws.Range('A348').Value = 'عوامل دسته a بیوتروریسم کدامند؟'
$ws.Range('A349').Value = 'کمبود هیدروکسیلاز با از دست دادن نمک با کدام یک از موارد زیر مشخص می‌شود؟  
الف) هیپوناترمی  
ب) هایپرکالمی  
ج) هیپوگلیسمی  
د) هیپوکلسمی'
$ws.Range('A350').Value = 'اولین شواهد کلسیفیکاسیون دندان‌های شیری قدامی تقریباً بین چه زمانی شروع می‌شود؟'
$ws.Range('A351').Value = 'ویژگی های سندرم هپاتورنال کدامند؟  
الف) سدیم ادرار کمتر از ۱۰ میلی‌اکیوالان در لیتر  
ب) هیستولوژی کلیه طبیعی  
ج) اختلال عملکرد کلیه حتی پس از بهبود کبد  
د) پروتئینوری کمتر از ۵۰۰ میلی‌گرم در روز'
$ws.Range('A352').Value = 'یک مرد 30 ساله با شکایت اصلی درد دائمی و تیرکشنده که با خوردن غذاهای تند و جویدن تشدید می‌شود، مراجعه کرده است. او همچنین از طعم "فلزی" بد و مقدار زیادی بزاق "خمیری" شکایت دارد. معاینه عمومی تب و افزایش ضربان قلب را نشان داد. معاینه داخل دهانی فرورفتگی‌های مشابه دهانه آتشفشان در لبه پاپی‌های بین دندانی که با یک لایه خاکستری کاذب پوشیده شده‌اند، در دندان‌های قدامی فک بالا مشاهده شد.

روش جراحی برای تخریب شدید بافت در چنین عفونتی چیست؟'
$ws.Range('A353').Value = 'فردی توسط یک شپش آلوده گزیده شده و دچار یک بیماری شده است. تشخیص احتمالی چیست؟'
$ws.Range('A354').Value = 'کدام یک از داروهای زیر در طول پرتودرمانی نیاز به تنظیم دوز دارد تا از سمیت پرتوی جلوگیری شود؟'
$ws.Range('A355').Value = 'مرحله‌بندی جراحی سرطان تخمدان تمام موارد زیر انجام می‌شود به جز؛'
$ws.Range('A356').Value = 'آنتاگونیسم بین استیل کولین و آتروپین:'
$ws.Range('A357').Value = 'کلیرانس یک ماده چقدر است، اگر غلظت آن در پلاسما 10 میلی‌گرم درصد، غلظت در ادرار 100 میلی‌گرم درصد و دفع ادرار 2 میلی‌لیتر در دقیقه باشد؟'
$ws.Range('A358').Value = 'سه‌گانه شارکو'
$ws.Range('A359').Value = 'آلوپورینول بر چه اساسی عمل می‌کند؟'
$ws.Range('A360').Value = 'کدام یک از موارد زیر در قربانیان تجاوز جنسی برای تعیین پارگی پرده بکارت و تشخیص جدید یا قدیمی بودن پارگی استفاده می‌شود؟'
$ws.Range('A361').Value = 'کدام عضله باعث دور شدن تارهای صوتی از هم می‌شود:'
$ws.Range('A362').Value = 'غنی‌ترین منبع کلسترول کدام است؟'
$ws.Range('A363').Value = 'در ارزیابی تغذیه‌ای یک بیمار جراحی، وضعیت پروتئین عضلانی توسط کدام یک از پارامترهای زیر نشان داده می‌شود؟'
$ws.Range('A364').Value = 'دلیل کنار گذاشتن روفکوکسیب چه بود؟'
$ws.Range('A365').Value = 'شایع ترین مسیر انتشار از محل اولیه که باعث سل رحم می شود کدام است؟'
$ws.Range('A366').Value = 'همه موارد زیر از نشانه‌های پیوند کبد هستند به جز -'
$ws.Range('A367').Value = 'یک کودک 3 ساله در صورتی به عنوان پنومونی طبقه‌بندی می‌شود که تعداد تنفس او باشد:'
$ws.Range('A368').Value = 'سندرم زلوگر یک اختلال مربوط به کدام است؟'
$ws.Range('A369').Value = 'کدام یک از موارد زیر از علائم بیماری انسدادی آترواسکلروتیک در دو شاخه شدن آئورت (سندرم لریش) نیست؟'
$ws.Range('A370').Value = 'ماده p متعلق به خانواده پپتیدهای تاکیکینین است. کدام یک از بافت های محیطی زیر حاوی ماده p هستند؟'
$ws.Range('A371').Value = 'آسپرین با کدام یک از موارد زیر مرتبط است؟'
$ws.Range('A372').Value = 'یک مرد 48 ساله مبتلا به ایدز با تب 38.7 درجه سانتیگراد (103 درجه فارنهایت)، سرفه مداوم و اسهال در بیمارستان بستری می‌شود. تعداد سلول‌های cd4 او کمتر از 500 در میلی‌لیتر است. برای بیمار آنتی‌بیوتیک‌های طیف وسیع شروع می‌شود. او همچنین اخیراً کاهش عملکرد شناختی را تجربه کرده است. خطر ابتلای او به کدام یک از نئوپلاسم‌های سیستم عصبی مرکزی زیر افزایش یافته است؟'
$ws.Range('A373').Value = 'عامل مستعد کننده برای سرطان معده به جز کدام یک از موارد زیر است؟'
$ws.Range('A374').Value = 'کدام ماده غذایی سرشار از آهن است؟'
$ws.Range('A375').Value = 'در مورد شدت درگیری اندام‌ها در دیپلژی اسپاستیک فلج مغزی کدام گزینه صحیح است؟'
$ws.Range('A376').Value = 'کدام یک از موارد زیر محصول متابولیک چرخه اوره نیست؟'
$ws.Range('A377').Value = 'مزایای نقدی در طرح esi شامل چه مواردی است؟'
$ws.Range('A378').Value = 'آنتاگونیست مواد افیونی که می‌توان به صورت خوراکی تجویز کرد –'
$ws.Range('A379').Value = 'کدام یک از افراد زیر ابتدا از اسید کربولیک به عنوان ضدعفونی کننده در جراحی استفاده کرد و به عنوان "پدر جراحی ضدعفونی کننده" شناخته می‌شود؟'
$ws.Range('A380').Value = 'کدام یک از موارد زیر در مورد اسفنکتر تحتانی مری نادرست است؟'
$ws.Range('A381').Value = 'حساسترین فیبر عصبی به هیپوکسی'
$ws.Range('A382').Value = 'آنتی‌ژن‌های یکسان یا نزدیک به هم که در گونه‌های بیولوژیکی مختلف وجود دارند، به عنوان چه شناخته می‌شوند؟'
$ws.Range('A383').Value = 'پاراکوزیس ویلیسی ویژگی کدام یک از موارد زیر است؟'
$ws.Range('A384').Value = 'آقای x به سمت آقای y شلیک کرد که حرکت کرد و فرار کرد و گلوله فقط بازوی او را خراش داد. خونریزی کمی وجود داشت و هیچ آسیب قابل توجه دیگری نبود. آقای x تحت کدام بخش از ipc قابل دستگیری است؟'
$ws.Range('A385').Value = 'نشانگر آسیب حاد کلیه به جز کدام یک؟'
$ws.Range('A386').Value = 'همه موارد زیر از عوامل خطر زایمان زودرس هستند به جز'
$ws.Range('A387').Value = 'مونو، یک مرد 30 ساله که الکلی مزمن است، با درد ناگهانی اپیگاستر که به پشت انتشار می‌یابد، مراجعه می‌کند. همه موارد زیر مشاهده می‌شوند به جز:'
$ws.Range('A388').Value = 'بیشترین خطر انتقال hiv به جنین در چه زمانی است؟'
$ws.Range('A389').Value = 'در مرور کلی بتا اکسیداسیون اسیدهای چرب، a نشان دهنده چیست؟'
$ws.Range('A390').Value = 'عضله سوپیناتور توسط کدام عصب عصب‌دهی می‌شود؟'
$ws.Range('A391').Value = 'خونرسانی به لوله استاش توسط کدام رگ انجام می‌شود؟'
$ws.Range('A392').Value = 'کدام گزینه در مورد سندرم مرگ ناگهانی نوزاد صحیح نیست؟'
$ws.Range('A393').Value = 'یک زن ۲۴ ساله در یک تصادف رانندگی قرار می‌گیرد و به اورژانس منتقل می‌شود، جایی که از قفسه سینه و ستون فقرات تحتانی او عکسبرداری می‌شود. بعداً مشخص می‌شود که او ۱۰ هفته باردار است. باید به او توصیه شود که:'
$ws.Range('A394').Value = 'یک بیمار مبتلا به زخم معده در بررسی آندوسکوپی، گاستریت آنترال مزمن نشان داد. کدام یک از رنگ‌های زیر قادر به رنگ‌آمیزی نمونه است؟'
$ws.Range('A395').Value = 'رایج‌ترین داروی مورد استفاده برای پیشگیری از میگرن کدام است؟'
$ws.Range('A396').Value = 'در مورد ورید لابه کدام گزینه صحیح است؟'
$ws.Range('A397').Value = 'غشای پلاسمایی یوکاریوتی از همه موارد زیر تشکیل شده است به جز:'
$ws.Range('A398').Value = 'شاخص درمانی دارو نشانگر چیست؟'
$ws.Range('A399').Value = 'ایستموس غده تیروئید در کدام سطح قرار دارد؟'
$ws.Range('A400').Value = 'عمل لاپاراسکوپی برای ریفلاکس گوارشی باید در بیمارانی با هر یک از شرایط زیر به جز یکی در نظر گرفته شود'
$ws.Range('A401').Value = 'کدام یک از سندرم های زیر با ضایعات قلبی همراه نیست؟'
$ws.Range('A402').Value = 'کدام گزینه برای پیشگیری از مالاریا استفاده نمی‌شود؟'
$ws.Range('A403').Value = 'رگ‌های خونی شبکیه از کدام قسمت تکوین می‌یابند؟'
$ws.Range('A404').Value = 'فردی به نام «x» در اثر تحریک، فرد دیگری به نام «y» را با چوبدستی میزند. این عمل منجر به تشکیل کبودی به اندازه 3 سانتیمتر در 3 سانتیمتر روی ساعد میشود. هیچ آسیب دیگری مشاهده نشده است. کدام یک از گزینههای زیر در مورد مجازات او صحیح است؟'
$ws.Range('A405').Value = 'اتساع آلوئول‌های دوردست در کدام نوع آمفیزم دیده می‌شود؟'
$ws.Range('A406').Value = 'پروتئین‌ها چگونه از دستگاه گوارش جذب می‌شوند؟'
$ws.Range('A407').Value = 'بهترین روش برای تشخیص دررفتگی تروماتیک عدسی چشم چیست؟'
$ws.Range('A408').Value = 'در مقایسه با چسب‌های اچ کامل، پرایمرهای خود اچ'
$ws.Range('A409').Value = 'یک زن ۲۱ ساله با الیگومنوره مراجعه می‌کند. در معاینه، رشد بیش از حد موهای صورت مشاهده می‌شود و آزمایشات آزمایشگاهی افزایش سطح تستوسترون آزاد سرم را نشان می‌دهند. سونوگرافی شکم انجام شده، تخمدان‌های طبیعی را نشان می‌دهد. محتمل‌ترین تشخیص چیست؟'
$ws.Range('A410').Value = 'شایع‌ترین علت سندرم موبیوس استفاده از کدام یک از داروهای زیر در دوران بارداری است؟'
$ws.Range('A411').Value = 'ساده ترین معیار پراکندگی'
$ws.Range('A412').Value = 'واکنش‌های لیکنوئید عمدتاً ناشی از چه چیزی هستند؟'
$ws.Range('A413').Value = 'کدام یک از داروهای زیر بیشترین احتمال جذب از مخاط معده را دارد؟'
$ws.Range('A414').Value = 'فردی که به مدت 24 سال برای دیابت نوع 2 تحت درمان بوده است، به عنوان بخشی از معاینه سالانه خود، آزمایش جمع‌آوری ادرار 24 ساعته انجام داده است. کاهش سطح کراتینین مشاهده شده است که به احتمال زیاد به دلیل کدام یک از موارد زیر است؟'
$ws.Range('A415').Value = 'بیمار با شکایت از سردرد و گرفتگی بینی به درمانگاه گوش و حلق و بینی مراجعه کرد. در سی‌تی اسکن، کدورت ناهمگن درگیر کننده چندین سینوس همراه با فرسایش استخوان (همانطور که در زیر نشان داده شده است) مشاهده شد. محتمل‌ترین تشخیص چیست؟'
$ws.Range('A416').Value = 'کدام یک از عبارات زیر در مورد هایپرپلازی گرهی کانونی (fnh) نادرست است؟'
$ws.Range('A417').Value = 'یک زن ۲۳ ساله برای پیشگیری از بارداری از قرص ضدبارداری استفاده می‌کند. کدام یک از عوامل زیر باعث می‌شود که یک متخصص سلامت، روش دیگری برای پیشگیری از بارداری به او توصیه کند؟'
$ws.Range('A418').Value = 'کدام یک از موارد زیر رتروپریتوان است؟'
$ws.Range('A419').Value = 'در هفته چهارم، جایی که آندودرم و اکتودرم در ناحیه سر و گردن به یکدیگر نزدیک می‌شوند کجاست؟'
$ws.Range('A420').Value = 'یک زن 25 ساله از کاهش شنوایی دوطرفه از 4 سال پیش شکایت دارد که در دوران بارداری تشدید شده است. نوع تمپانوگرام در این حالت چیست؟'
$ws.Range('A421').Value = 'عبارت نادرست در مورد لیتیوم کدام است؟'
$ws.Range('A422').Value = 'رشد جذابیت در دختران جوان نسبت به پدرشان:'
$ws.Range('A423').Value = 'بیوپسی با سوزن از عصب سورال منجر به تشکیل هماتوم شده است. کدام یک از رگ‌های زیر که در مجاورت نزدیک این عصب قرار دارد، به طور تصادفی آسیب دیده است؟'
$ws.Range('A424').Value = 'کدام یک از موارد زیر از ویژگی‌های شبه‌کم‌کاری پاراتیروئید نیست؟'
$ws.Range('A425').Value = 'کدام میکروارگانیسم به عنوان سلاح در تروریسم بیولوژیک استفاده می‌شود؟'
$ws.Range('A426').Value = 'در مسمومیت با جیوه، رفلکس قهوه‌ای از کدام ناحیه ناشی می‌شود؟ ap 08'
$ws.Range('A427').Value = 'داروی مورد استفاده برای تسکین درد در بخش بیماران سرپایی (opd) کدام است؟'
$ws.Range('A428').Value = 'یک مرد ۳۵ ساله با سندرم داون به دلیل لوسمی لنفوبلاستیک حاد فوت می‌کند. بررسی ماکروسکوپی مغز بیمار در کالبدشکافی، میکروسفالی خفیف و عدم تکامل کافی چین‌های گیجگاهی فوقانی را نشان می‌دهد. بررسی هیستولوژیک به احتمال زیاد کدام یک از تغییرات نوروپاتولوژیک زیر را نشان خواهد داد؟'
$ws.Range('A429').Value = 'شایع‌ترین عارضه جانبی جراحی کاشت لنز چیست؟'
$ws.Range('A430').Value = 'در سندرم شریان نخاعی قدامی کدام یک از موارد زیر سالم می‌ماند؟'
$ws.Range('A431').Value = 'کدام یک از بی‌حس‌کننده‌های موضعی زیر متعلق به گروه استر است؟'
$ws.Range('A432').Value = 'یک کودک 9 ساله با شکایت از تب بالا، استفراغ و یک بار تشنج به درمانگاه مراجعه کرده است. آزمایش csf انجام شد و رنگ‌آمیزی گرم کشت، یافته زیر را نشان داد. عامل احتمالی ایجادکننده چیست؟'
$ws.Range('A433').Value = 'روش غربالگری انتخابی در منطقه‌ای که شیوع جذام 1 در 1000 است، کدام است؟'
$ws.Range('A434').Value = 'تست ویتاکر برای تشخیص کدام یک از شرایط زیر انجام می‌شود؟'
$ws.Range('A435').Value = 'کدام یک از داروهای زیر برای ترک سیگار استفاده می‌شوند؟'
$ws.Range('A436').Value = 'اثر پروژسترون چیست؟'
$ws.Range('A437').Value = 'اظهارات شخص در حال مرگ تحت کدام بخش قرار می‌گیرد؟'
$ws.Range('A438').Value = 'تست بانی برای چه موردی انجام می‌شود؟'
$ws.Range('A439').Value = 'تمامی موارد زیر از اولویت‌های فوری برنامه "بینایی ۲۰۲۰: حق بینایی" سازمان جهانی بهداشت هستند به جز:'
$ws.Range('A440').Value = 'سندرم ری توسط همه موارد زیر به جز کدام یک ایجاد می‌شود؟'
$ws.Range('A441').Value = 'یک زن ۷۲ ساله در یک ویزیت معمولی کلینیک، دارای ضربان نامنظم تشخیص داده می‌شود. او در حالت استراحت یا فعالیت هیچ علامت جدیدی را تجربه نمی‌کند. سابقه پزشکی او شامل فشار خون بالا، استئوآرتریت و دیس‌لیپیدمی است. در معاینه فیزیکی، فشار خون ۱۳۵/۸۵ میلی‌متر جیوه، ضربان قلب تقریباً ۷۲ در دقیقه و نامنظم است. صداهای قلب نشان‌دهنده s1 نامنظم و s2 طبیعی بدون سوفل قابل شنیدن است. فشار ورید ژوگولار طبیعی و ریه‌ها صاف هستند. در نوار قلب، موج p وجود ندارد و فاصله rr نامنظم با سرعت ۷۰ در دقیقه است. (شکل زیر را ببینید) در نوار قلب قبلی او که مربوط به ۴ سال پیش است، ریتم سینوسی داشت. کدام یک از گزینه‌های زیر مناسب‌ترین اقدام بعدی در مدیریت این بیمار است؟'
$ws.Range('A442').Value = 'یک بیمار 40 ساله با یک کلیه منفرد، با یک توده برون‌رویی منفرد به اندازه 4 سانتیمتر در قطب تحتانی آن مراجعه می‌کند. کدام یک از گزینه‌های زیر بهترین گزینه مدیریتی توصیه شده است؟'
$ws.Range('A443').Value = 'منع مصرف مطلق قرص ضدبارداری خوراکی کدام است؟'
$ws.Range('A444').Value = 'اکسید نیتریک توسط چه چیزی تولید می‌شود؟'
$ws.Range('A445').Value = 'همه موارد زیر ضایعات تاولی-وزیکولی هستند به جز-'
$ws.Range('A446').Value = 'تب خونریزی دهنده دنگی ناشی از چیست؟'
$ws.Range('A447').Value = 'آویزان کردن باعث آسیب زیاد به کدام قسمت می‌شود؟'
$ws.Range('A448').Value = 'آبسه سیته‌لی چیست؟'
$ws.Range('A449').Value = 'اعمال نیروی منفرد در کدام نقطه از دندان باعث انتقال کامل دندان می‌شود؟'
$ws.Range('A450').Value = 'نشانگر تومور سلول‌های گرانولوزا:'
$ws.Range('A451').Value = 'کدام سلول‌ها در مسیر بینایی به‌صورت خالص دپولاریزه می‌شوند؟'
$ws.Range('A452').Value = 'قلب در حالت استراحت از چه چیزی استفاده می‌کند؟'
$ws.Range('A453').Value = 'خطر انتقال جفتی توکسوپلاسما گوندی زمانی که مادر آلوده است، در چه زمانی بیشترین است؟'
$ws.Range('A454').Value = 'عامل مستعد کننده برای سرطان مری همه موارد زیر به جز کدام است؟'
$ws.Range('A455').Value = 'ارکیدوپکسی برای بیضه نزول نکرده کامل پس از چه سنی انجام می‌شود؟ سپتامبر 2011'
$ws.Range('A456').Value = 'مرحله انقباض خانواده در چه زمانی است؟'
$ws.Range('A457').Value = 'یک کودک ۱۰ ساله در حال انجام عمل جراحی انحراف چشم است. او ناگهان دچار افزایش ضربان قلب، آریتمی، تب بالا، اسیدوز متابولیک و تنفسی در گازهای خون شریانی و افزایش co2 پایان بازدمی می‌شود. کدام یک از موارد زیر اولین انتخاب در مدیریت این وضعیت است؟'
$ws.Range('A458').Value = 'در ارزیابی معمول آزمایشگاهی، بیمار از نظر hbs ag مثبت است. سایر آزمایش‌های سرولوژیکی برای هپاتیت قابل توجه نیستند. بیمار از نظر بالینی بدون علامت است و آنزیم‌های کبدی در محدوده طبیعی قرار دارند. کدام یک از گزینه‌های زیر بهترین توصیف برای تشخیص این بیمار است؟'
$ws.Range('A459').Value = 'تمام موارد زیر در مورد کلامیدیا پسیتاسی صحیح هستند به جز:'
$ws.Range('A460').Value = 'کدام یک تحریک کننده است؟'
$ws.Range('A461').Value = 'ضایعات دهانی در حالت نقص ایمنی'
$ws.Range('A462').Value = 'یک نوزاد با تشنجهای مقاوم به درمان و صدای پیوسته از طریق فونتانل قدامی آورده شده است. سی‌تی‌اسکن یک ضایعه خط میانی با هایپواکویستی و بطن‌های جانبی گشاد شده را نشان می‌دهد. محتمل‌ترین تشخیص چیست؟'
$ws.Range('A463').Value = 'هایپرونتیلاسیون در ارتفاعات بالا به دلیل چیست؟'
$ws.Range('A464').Value = 'جهش لیدن شامل کدام مورد است؟'
$ws.Range('A465').Value = 'زاویه کاب برای اندازه گیری میزان کدام مورد استفاده می شود؟'
$ws.Range('A466').Value = 'فلج ارب مربوط به کدام بخش است؟'
$ws.Range('A467').Value = 'تمام موارد زیر در مورد بیماری کلیه پلی‌کیستیک بزرگسالان صحیح است، به جز:'
$ws.Range('A468').Value = 'دندان مصنوعی جیفی نوعی از کدام مورد است؟'
$ws.Range('A469').Value = 'مرد 58 ساله با شروع حاد واریکوسل در سمت چپ، محتمل‌ترین علت چیست؟'
$ws.Range('A470').Value = 'اکسی توسین همه موارد زیر را باعث می‌شود به جز'
$ws.Range('A471').Value = 'کاریوتایپینگ در کدام مرحله از چرخه سلولی انجام می‌شود؟'
$ws.Range('A472').Value = 'سازمان جهانی بهداشت اعلام کرد که آبله در چه تاریخی ریشه‌کن شده است؟'
$ws.Range('A473').Value = 'همه ساختارهای زیر از فاشیای کلاویپکتورال عبور میکنند، به جز:'
$ws.Range('A474').Value = 'سل لوله های فالوپ -'
$ws.Range('A475').Value = 'در یک کودک با نارسایی حاد کبدی، مهمترین آزمایش بیوشیمیایی غیرطبیعی سرم که نشان‌دهنده پیش‌آگهی ضعیف است، کدام است؟'
$ws.Range('A476').Value = 'موارد زیر با فئوکروموسیتوما سازگار هستند به جز -'
$ws.Range('A477').Value = 'فیبرهای پروپریوسپتیو توسط کدام عصب جمجمه‌ای منتقل نمی‌شوند؟'
$ws.Range('A478').Value = 'انوفتالموس ناشی از شکستگی کدام قسمت است؟'
$ws.Range('A479').Value = 'رسوبات کراتیک در کدام لایه قرنیه قرار دارند؟'
$ws.Range('A480').Value = 'مقدار کل آهن مورد نیاز جنین در طول دوران بارداری'
$ws.Range('A481').Value = 'بهترین روش تشخیصی برای سنگ های حالب چیست؟'
$ws.Range('A482').Value = 'بند ناف دارای ____________ است'
$ws.Range('A483').Value = 'تمام موارد زیر در مورد پشه aedes صحیح است، به جز:'
$ws.Range('A484').Value = 'مخرج نرخ مرگ و میر مادران چیست؟ سپتامبر 2005، 2010، مارس 2013 (d, e, h)'
$ws.Range('A485').Value = 'استخوانچه های گوش میانی مسئول کدام یک از موارد زیر هستند؟'
$ws.Range('A486').Value = 'واسکولیت هایپرسنسیتیویتی در کدام یک از موارد زیر دیده می‌شود؟'
$ws.Range('A487').Value = 'اصل "برای ایجاد تغییر در پاسخ، محرک باید تغییر کند" در کدام مورد اعمال می‌شود؟'
$ws.Range('A488').Value = 'یک کارگر چندمنظوره زن باید بتواند همه موارد زیر را تشخیص دهد به جز:'
$ws.Range('A489').Value = 'یک مرد 50 ساله سابقه حملات مکرر کولیک کلیوی با سنگ‌های کلیوی حاوی کلسیم بالا دارد. مفیدترین دیورتیک در درمان سنگ‌های کلسیمی مکرر کدام است؟'
$ws.Range('A490').Value = 'بیماری پاژه همچنین به عنوان چه چیزی شناخته می‌شود؟'
$ws.Range('A491').Value = 'لکه‌های کوپلیک نشانه‌ی پاتوگنومونیک کدام عفونت است؟'
$ws.Range('A492').Value = 'شایع‌ترین علت استریدور متناوب در یک کودک ۱۰ روزه چیست؟'
$ws.Range('A493').Value = 'کدام دیورتیک بدون نیاز به دسترسی به لومن توبول کلیوی عمل می‌کند؟'
$ws.Range('A494').Value = 'در مورد فرمولاسیون لیپیدی یا لیپوزومی آمفوتریسین b کدام یک از عبارات زیر صحیح است؟'
$ws.Range('A495').Value = 'حداقل غلظت آلوئولی سِووفلوران بر حسب درصد چقدر است؟'
$ws.Range('A496').Value = 'یک پاکت پریودنتال به عمق 6 میلی‌متر که در قسمت کرونال اپی‌تلیوم اتصال قرار دارد، چیست؟'
$ws.Range('A497').Value = 'کدام یک از نشانگرهای زیر برای تشخیص رابدومیوسارکوم استفاده می‌شود؟'
$ws.Range('A498').Value = 'تصعید تبدیل چه چیزی است؟'
$ws.Range('A499').Value = 'استراتژی safe برای تراخم شامل همه موارد زیر به جز کدام است؟'
$ws.Range('A500').Value = 'چرا انتقال اکسیژن توسط هموگلوبین در اثر کاهش فعالیت گلیکولیتیک مختل می‌شود؟'
$ws.Range('A501').Value = 'داروی انتخابی در تنگی حنجره کدام است؟'
$ws.Range('A502').Value = 'سندرم ریه کوچک شده در کدام یک از موارد زیر دیده می‌شود؟'
$ws.Range('A503').Value = 'نظریه کنترل دروازه‌ای درد مربوط به کدام مورد است؟'
$ws.Range('A504').Value = 'واکسن‌های القاکننده ctl، واکسن ویروس آدنو-مرتبط نوترکیب و واکسن اصلاح‌شده آنکارا برای چه بیماری‌هایی در حال توسعه هستند؟'
$ws.Range('A505').Value = 'یک خانم 35 ساله به مدت 4 ماه گذشته قاعدگی نداشته است. سطح سرمی fsh و lh او بالا و استرادیول او پایین است. علت احتمالی چیست؟'
$ws.Range('A506').Value = 'اسمگما توسط کدام غده ترشح می‌شود؟'
$ws.Range('A507').Value = 'لکه‌های میلاری در ریه در کدام مورد دیده می‌شود؟'
$ws.Range('A508').Value = 'شایع‌ترین علت نورون حرکتی تحتانی (lmn) فلج عصب صورت چیست؟'
$ws.Range('A509').Value = 'غدد لنفاوی مثلث خلفی جزو کدام سطح هستند؟'
$ws.Range('A510').Value = 'ناهنجاری مادرزادی ناشی از درمان با لیتیوم کدام است؟'
$ws.Range('A511').Value = 'تومور وارین چیست؟'
$ws.Range('A512').Value = 'ترشح هورمون رشد توسط همه موارد زیر کاهش می‌یابد به جز'
$ws.Range('A513').Value = 'پدیده پروزون به دلیل چیست؟  
الف) آنتی‌ژن اضافی  
ب) آنتی‌بادی اضافی  
ج) واکنش بیش‌ازحد ایمنی  
د) سطوح نامتناسب آنتی‌ژن-آنتی‌بادی'
$ws.Range('A514').Value = 'مدت زمان ضدعفونی کردن قالب آلژیناتی، در صورت استفاده از روش غوطه‌وری، نباید از چه مدت زمانی بیشتر شود؟'
$ws.Range('A515').Value = 'کدام دارو به عنوان داروی پیش بیهوشی استفاده می‌شود؟'
$ws.Range('A516').Value = 'نوزاد با وزن بسیار کم هنگام تولد (vlbw) نوزادی است که وزن هنگام تولد آن کمتر از ____ کیلوگرم باشد.'
$ws.Range('A517').Value = 'نمره gcs در یک بیمار مرده چقدر است؟'
